# Semana 19 de 2025: add column V ("19") to the weekly IRA extract sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell V1 is a text label (matching "1".."18" in D1:U1), entered
# with a leading apostrophe so it is stored as text, not a number.
$ws.Range("V1").Formula = "'19"

# Weekly case counts for week 19 (column V), row by row.
$ws.Range("V2").Value = 57
$ws.Range("V3").Value = 49
$ws.Range("V4").Value = 0
$ws.Range("V5").Value = 5
$ws.Range("V6").Value = 52
$ws.Range("V7").Value = 25
$ws.Range("V8").Value = 33
$ws.Range("V9").Value = 3
$ws.Range("V10").Value = 5
$ws.Range("V11").Value = 1
$ws.Range("V13").Value = 2
$ws.Range("V14").Value = 3
$ws.Range("V15").Value = 2
$ws.Range("V17").Value = 1
$ws.Range("V20").Value = 1
$ws.Range("V23").Value = 55
$ws.Range("V24").Value = 1
$ws.Range("V26").Value = 178
$ws.Range("V27").Value = 0
$ws.Range("V28").Value = 22
$ws.Range("V29").Value = 5
$ws.Range("V30").Value = 5
$ws.Range("V32").Value = 37
$ws.Range("V33").Value = 4
$ws.Range("V34").Value = 5
$ws.Range("V35").Value = 75
$ws.Range("V37").Value = 12
$ws.Range("V38").Value = 48
$ws.Range("V39").Value = 25
$ws.Range("V40").Value = 192
$ws.Range("V41").Value = 91
$ws.Range("V42").Value = 153
$ws.Range("V43").Value = 5
$ws.Range("V44").Value = 96
$ws.Range("V45").Value = 3
$ws.Range("V46").Value = 0
$ws.Range("V47").Value = 4
$ws.Range("V48").Value = 1
$ws.Range("V49").Value = 41
$ws.Range("V50").Value = 0
$ws.Range("V51").Value = 0
$ws.Range("V52").Value = 4
$ws.Range("V53").Value = 14
$ws.Range("V54").Value = 46
